# Quarterly income-statement update: a new quarter ("1401/12" -> "1402/01/29
# (8)") is appended, every existing quarter column shifts one column to the
# left (D<-E<-F...<-M), and a handful of derived figures are recomputed
# (read_price algorithm change) rather than simply carried over by the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: quarter labels (D8:M8) ---------------------------------------
$row8 = @(
    "فصل سوم منتهی به 1399/09",
    "فصل چهارم منتهی به 1399/12",
    "فصل اول منتهی به 1400/03",
    "فصل دوم منتهی به 1400/06",
    "فصل سوم منتهی به 1400/09",
    "فصل چهارم منتهی به 1400/12",
    "فصل اول منتهی به 1401/03",
    "فصل دوم منتهی به 1401/06",
    "فصل سوم منتهی به 1401/09",
    "فصل چهارم منتهی به 1401/12"
)
for ($i = 0; $i -lt $row8.Length; $i++) {
    $ws.Cells.Item(8, 4 + $i).Value = $row8[$i]
}

# --- Row 9: publish dates (D9:M9) -----------------------------------------
$row9 = @(
    "1400-10-30 (2)",
    "1401-04-08 (9)",
    "1401-04-30 (2)",
    "1401-09-15 (4)",
    "1401-10-28 (2)",
    "1402-01-29 (8)",
    "1401-04-30",
    "1401-09-15 (2)",
    "1401-10-28",
    "1402-01-29"
)
for ($i = 0; $i -lt $row9.Length; $i++) {
    $ws.Cells.Item(9, 4 + $i).Value = $row9[$i]
}

# --- Data rows 11-27 (D:M), shifted left by one quarter --------------------
# Each entry is the new D..M content for that row (10 values).
$dataRows = @{
    11 = @(301888, 489511, 418187, 449733, 423808, 437028, 422195, 266458, 295513, 222485)
    12 = @(-95460, -159217, -137174, -151552, -218744, -118179, -176732, -174302, -188631, -104628)
    13 = @(206428, 269492, 281013, 298181, 205064, 318848, 245463, 92156, 106882, 117857)
    14 = @(-764, -10869, -1281, -11374, -5640, -9803, -6613, -6879, -6553, -7913)
    15 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    16 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    17 = @(205664, 258623, 279732, 286807, 199424, 309045, 238850, 85277, 100329, 109944)
    18 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    19 = @(1688, 15590, 11058, 9179, 11555, 10585, 12272, 14152, 28049, 11576)
    20 = @(207353, 335015, 290790, 295987, 210979, 319630, 251123, 99429, 128378, 121520)
    21 = @(-18447, -22195, -33055, -32323, -23388, -3339, -16515, -12028, -1440, 19985)
    22 = @(188905, 312819, 257735, 263664, 187591, 316291, 234608, 87401, 126938, 141505)
    23 = @("-", "-", "-", "-", "-", "-", "-", "-", "-", "-")
    24 = @(188905, 312819, 257735, 263664, 187591, 316291, 234608, 87401, 126938, 141505)
    25 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
    26 = @(204224, 227232, 237892, 212634, 417962, 432040, 404972, 383675, 342516, 423513)
    27 = @(0, 0, 0, 0, 0, 0, 0, 0, 0, 0)
}

foreach ($r in $dataRows.Keys) {
    $vals = $dataRows[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, 4 + $i).Value = $vals[$i]
    }
}

# --- Column widths also shift left by one quarter (D:M) --------------------
# stored xlsx <col> width = ColumnWidth + 0.8333333333333321 on this engine,
# so back the offset out to land on the exact target widths (29 / 31).
$narrow = 29 - 0.8333333333333321
$wide   = 31 - 0.8333333333333321
$ws.Columns.Item(5).ColumnWidth  = $wide    # E: 29 -> 31
$ws.Columns.Item(6).ColumnWidth  = $narrow  # F: 31 -> 29
$ws.Columns.Item(9).ColumnWidth  = $wide    # I: 29 -> 31
$ws.Columns.Item(10).ColumnWidth = $narrow  # J: 31 -> 29
$ws.Columns.Item(13).ColumnWidth = $wide    # M: 29 -> 31
